$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.1
$ws.Range("G2").Value = 1.11
$ws.Range("H2").Value = 17.5
$ws.Range("I2").Value = 980
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 110
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 3.45
$ws.Range("Q2").Value = 1.37
$ws.Range("R2").Value = 1.46
$ws.Range("S2").Value = 2.84
$ws.Range("T2").Value = 2.42
$ws.Range("U2").Value = 1.59
$ws.Range("V2").Value = 1.02
$ws.Range("W2").Value = 10
$ws.Range("AC2").Value = 1000
$ws.Range("AF2").Value = 3.75
$ws.Range("AG2").Value = 980
$ws.Range("AH2").Value = 980
$ws.Range("AI2").Value = 220
$ws.Range("AJ2").Value = 10.5
$ws.Range("AK2").Value = 12
$ws.Range("AN2").Value = 980
$ws.Range("F3").Value = 1.37
$ws.Range("G3").Value = 1.41
$ws.Range("H3").Value = 12.5
$ws.Range("I3").Value = 15.5
$ws.Range("J3").Value = 4.6
$ws.Range("O3").Value = 1.43
$ws.Range("Q3").Value = 2.3
$ws.Range("T3").Value = 2.74
$ws.Range("U3").Value = 1.54
$ws.Range("V3").Value = 1.07
$ws.Range("W3").Value = 3.4
$ws.Range("Y3").Value = 29
$ws.Range("AD3").Value = 190
$ws.Range("AH3").Value = 50
$ws.Range("F4").Value = 1.36
$ws.Range("G4").Value = 1.4
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 5.6
$ws.Range("L4").Value = 1.23
$ws.Range("N4").Value = 6.4
$ws.Range("O4").Value = 1.12
$ws.Range("P4").Value = 2.96
$ws.Range("Q4").Value = 1.41
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 2.02
$ws.Range("H5").Value = 3.05
$ws.Range("I5").Value = 3.85
$ws.Range("N5").Value = 3.45
$ws.Range("O5").Value = 1.28
$ws.Range("P5").Value = 1.89
$ws.Range("S5").Value = 2.88
$ws.Range("U5").Value = 2.08
$ws.Range("W5").Value = 1.66
$ws.Range("Y5").Value = 15
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 65
$ws.Range("AE5").Value = 42
$ws.Range("AI5").Value = 50
$ws.Range("AO5").Value = 40
$ws.Range("H6").Value = 1.47
$ws.Range("L6").Value = 1.38
$ws.Range("N6").Value = 4
$ws.Range("O6").Value = 1.31
$ws.Range("Q6").Value = 1.92
$ws.Range("S6").Value = 3.35
$ws.Range("U6").Value = 1.84
$ws.Range("X6").Value = 17
$ws.Range("Z6").Value = 8
$ws.Range("AD6").Value = 9.800000000000001
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 300
$ws.Range("AM6").Value = 160
$ws.Range("AO6").Value = 8
$ws.Range("J7").Value = 3.35
$ws.Range("O7").Value = 1.3
$ws.Range("Q7").Value = 1.76
$ws.Range("AH7").Value = 1000
$ws.Range("G8").Value = 1.82
$ws.Range("H8").Value = 5.1
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 3.95
$ws.Range("N8").Value = 3.6
$ws.Range("P8").Value = 1.87
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.34
$ws.Range("S8").Value = 3.35
$ws.Range("U8").Value = 1.91
$ws.Range("W8").Value = 2.2
$ws.Range("Y8").Value = 20
$ws.Range("AA8").Value = 150
$ws.Range("AB8").Value = 8.6
$ws.Range("AC8").Value = 9
$ws.Range("AI8").Value = 85
$ws.Range("AN8").Value = 12.5
